$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1:G1").EntireColumn.Delete()
$ws.Range("C7").Select()
$ws.Protect("password")
